$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Google.com" entry (row 2) to lowercase "google.com"
$ws.Range("A2").Value = "google.com"

# Move selection to A2, matching the saved cursor position
$ws.Range("A2").Select() | Out-Null
